$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative xml diff for this commit.
$updates = @(
    @{ Cell = "D2"; Value = "30.075.81" },
    @{ Cell = "E2"; Value = "  -0.82%  " },
    @{ Cell = "D3"; Value = "1.906.26" },
    @{ Cell = "E3"; Value = "  -1.45%  " },
    @{ Cell = "E4"; Value = "  -0.12%  " },
    @{ Cell = "D5"; Value = "0.7440" },
    @{ Cell = "E5"; Value = "  -0.75%  " },
    @{ Cell = "D6"; Value = "244.03" },
    @{ Cell = "E6"; Value = "  +0.55%  " },
    @{ Cell = "D7"; Value = "1.000" },
    @{ Cell = "E7"; Value = "  -0.12%  " },
    @{ Cell = "D8"; Value = "0.3095" },
    @{ Cell = "E8"; Value = "  -2.59%  " },
    @{ Cell = "D9"; Value = "26.45" },
    @{ Cell = "E9"; Value = "  -5.36%  " },
    @{ Cell = "D10"; Value = "0.06978" },
    @{ Cell = "E10"; Value = "  -3.31%  " },
    @{ Cell = "D11"; Value = "0.08086" },
    @{ Cell = "E11"; Value = "  +0.59%  " },
    @{ Cell = "D12"; Value = "0.7682" },
    @{ Cell = "E12"; Value = "  -1.41%  " },
    @{ Cell = "D13"; Value = "1.930.00" },
    @{ Cell = "E13"; Value = "  -0.23%  " },
    @{ Cell = "D14"; Value = "5.301" },
    @{ Cell = "D15"; Value = "92.13" },
    @{ Cell = "E15"; Value = "  -0.91%  " },
    @{ Cell = "D16"; Value = "14.22" },
    @{ Cell = "E16"; Value = "  -1.90%  " },
    @{ Cell = "D17"; Value = "30.073.77" },
    @{ Cell = "E17"; Value = "  -0.89%  " },
    @{ Cell = "D18"; Value = "6.060" },
    @{ Cell = "E18"; Value = "  -0.91%  " },
    @{ Cell = "D19"; Value = "0.000007827" },
    @{ Cell = "E19"; Value = "  -2.40%  " },
    @{ Cell = "D20"; Value = "239.74" },
    @{ Cell = "E20"; Value = "  -4.79%  " },
    @{ Cell = "B21"; Value = "WrappedliquidstakedEther2.0" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" },
    @{ Cell = "D21"; Value = "2.166.04" },
    @{ Cell = "E21"; Value = "  -0.98%  " },
    @{ Cell = "B22"; Value = "Dai" },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai" },
    @{ Cell = "D22"; Value = "1.001" },
    @{ Cell = "E22"; Value = "  +0.06%  " },
    @{ Cell = "E23"; Value = "  -0.04%  " },
    @{ Cell = "D24"; Value = "7.149" },
    @{ Cell = "E24"; Value = "  +6.94%  " },
    @{ Cell = "D25"; Value = "9.369" },
    @{ Cell = "E25"; Value = "  -1.94%  " },
    @{ Cell = "D26"; Value = "167.01" },
    @{ Cell = "E26"; Value = "  +1.30%  " },
    @{ Cell = "E27"; Value = "  -0.45%  " },
    @{ Cell = "D28"; Value = "0.1270" },
    @{ Cell = "E28"; Value = "  -2.48%  " },
    @{ Cell = "D29"; Value = "2.052" },
    @{ Cell = "E29"; Value = "  -6.57%  " },
    @{ Cell = "E30"; Value = "  -1.98%  " },
    @{ Cell = "D31"; Value = "1.543" },
    @{ Cell = "E31"; Value = "  -0.04%  " },
    @{ Cell = "D32"; Value = "4.335" },
    @{ Cell = "E32"; Value = "  -2.09%  " },
    @{ Cell = "D33"; Value = "4.078" },
    @{ Cell = "E33"; Value = "  -1.66%  " },
    @{ Cell = "D34"; Value = "0.05220" },
    @{ Cell = "E34"; Value = "  -1.26%  " },
    @{ Cell = "E35"; Value = "  -2.40%  " },
    @{ Cell = "D36"; Value = "0.7472" },
    @{ Cell = "E36"; Value = "  -1.13%  " },
    @{ Cell = "E37"; Value = "  -2.50%  " },
    @{ Cell = "E38"; Value = "  +0.45%  " },
    @{ Cell = "E39"; Value = "  -0.02%  " },
    @{ Cell = "D40"; Value = "6.324" },
    @{ Cell = "E40"; Value = "  -2.62%  " },
    @{ Cell = "D41"; Value = "0.4483" },
    @{ Cell = "E41"; Value = "  -0.60%  " },
    @{ Cell = "D42"; Value = "74.17" },
    @{ Cell = "E42"; Value = "  -5.75%  " },
    @{ Cell = "D43"; Value = "1.973" },
    @{ Cell = "E43"; Value = "  -0.43%  " },
    @{ Cell = "E44"; Value = "  +0.03%  " },
    @{ Cell = "D45"; Value = "0.8399" },
    @{ Cell = "E45"; Value = "  -0.03%  " },
    @{ Cell = "D46"; Value = "7.721" },
    @{ Cell = "E46"; Value = "  +0.30%  " },
    @{ Cell = "D47"; Value = "101.72" },
    @{ Cell = "E47"; Value = "  +0.16%  " },
    @{ Cell = "D48"; Value = "9.894" },
    @{ Cell = "E48"; Value = "  -1.00%  " },
    @{ Cell = "D49"; Value = "2.070.24" },
    @{ Cell = "E49"; Value = "  -1.06%  " },
    @{ Cell = "D50"; Value = "36.62" },
    @{ Cell = "E50"; Value = "  -2.63%  " },
    @{ Cell = "D51"; Value = "0.1180" },
    @{ Cell = "E51"; Value = "  -4.14%  " }
)

foreach ($u in $updates) {
    # Leading apostrophe forces Excel to treat the assignment as literal text even
    # when the string looks numeric/date-like (e.g. "1.000", "0.7440", "5.301");
    # these price/volume cells are stored as text in the workbook, not numbers.
    $ws.Range($u.Cell).Value = "'" + $u.Value
    # Reset to the workbook-default style so the forced-text number format used
    # above does not linger as a visible/applied style on the cell.
    $ws.Range($u.Cell).Style = "Normal"
}
